$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Col 6" header in J5 (cell becomes blank / no value)
$ws.Range("J5").Value = $null

# G6 was "admin" -> now "admin2"
$ws.Range("G6").Value = "admin2"

# G7 was "admin2" -> now "admin1"
$ws.Range("G7").Value = "admin1"

# Update the active selection to J5 (was J8)
$ws.Range("J5").Select()
